$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "68.187.08"
Set-TextValue "E2" "  -0.06%  "
Set-TextValue "D3" "3.623.58"
Set-TextValue "E3" "  -1.11%  "
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "588.16"
Set-TextValue "E5" "  -1.30%  "
Set-TextValue "D6" "194.31"
Set-TextValue "E6" "  +1.11%  "
Set-TextValue "D7" "3.619.63"
Set-TextValue "E7" "  -0.97%  "
Set-TextValue "D8" "0.621"
Set-TextValue "E8" "  -0.09%  "
Set-TextValue "E9" "  -0.06%  "
Set-TextValue "D10" "0.682"
Set-TextValue "E10" "  -2.28%  "
Set-TextValue "D11" "0.153"
Set-TextValue "E11" "  -0.14%  "
Set-TextValue "D12" "55.69"
Set-TextValue "E12" "  -2.47%  "
Set-TextValue "D13" "0.0000292"
Set-TextValue "E13" "  +7.52%  "
Set-TextValue "D14" "10.03"
Set-TextValue "E14" "  -2.06%  "
Set-TextValue "D15" "4.191.66"
Set-TextValue "E15" "  -1.46%  "
Set-TextValue "D16" "3.620.44"
Set-TextValue "E17" "  -0.37%  "
Set-TextValue "D18" "12.58"
Set-TextValue "E18" "  -0.16%  "
Set-TextValue "D19" "68.021.99"
Set-TextValue "E19" "  -0.03%  "
Set-TextValue "D20" "18.58"
Set-TextValue "E20" "  -1.53%  "
Set-TextValue "E21" "  -2.19%  "
Set-TextValue "D22" "405.71"
Set-TextValue "E22" "  +0.19%  "
Set-TextValue "D23" "13.50"
Set-TextValue "E23" "  +24.11%  "
Set-TextValue "D24" "4.26"
Set-TextValue "E24" "  -3.30%  "
Set-TextValue "D25" "86.20"
Set-TextValue "E25" "  -2.42%  "
Set-TextValue "E26" "  +0.58%  "
Set-TextValue "D27" "12.64"
Set-TextValue "E27" "  +0.53%  "
Set-TextValue "D28" "3.93"
Set-TextValue "E28" "  +5.39%  "
Set-TextValue "E29" "  +0.90%  "
Set-TextValue "D30" "8.27"
Set-TextValue "E30" "  +15.38%  "
Set-TextValue "D31" "9.21"
Set-TextValue "E31" "  -1.40%  "
Set-TextValue "D32" "31.72"
Set-TextValue "E32" "  -0.84%  "
Set-TextValue "D33" "680.24"
Set-TextValue "E33" "  +12.19%  "
Set-TextValue "D34" "12.28"
Set-TextValue "E34" "  +0.06%  "
Set-TextValue "E35" "  +1.73%  "
Set-TextValue "D36" "64.64"
Set-TextValue "E36" "  -3.68%  "
Set-TextValue "D37" "42.60"
Set-TextValue "E37" "  -3.09%  "
Set-TextValue "D38" "0.425"
Set-TextValue "E38" "  +8.34%  "
Set-TextValue "E39" "  +0.04%  "
Set-TextValue "D40" "0.0₃0794"
Set-TextValue "E40" "  +2.53%  "
Set-TextValue "D41" "2.98"
Set-TextValue "E41" "  +18.55%  "
Set-TextValue "E42" "  +8.40%  "
Set-TextValue "D43" "3.208.07"
Set-TextValue "E43" "  +15.50%  "
Set-TextValue "D44" "0.135"
Set-TextValue "E44" "  -1.06%  "
Set-TextValue "D45" "0.999"
Set-TextValue "E45" "  -0.21%  "
Set-TextValue "E46" "  -0.97%  "
Set-TextValue "E47" "  -2.37%  "

# Rows 48-50 reorder: THORChain / Monero / ApeXProtocol shuffle order & values
Set-TextValue "B48" "ApeXProtocol"
Set-TextValue "C48" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D48" "3.13"
Set-TextValue "E48" "  -1.80%  "

Set-TextValue "B49" "THORChain"
Set-TextValue "C49" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D49" "8.83"
Set-TextValue "E49" "  -1.08%  "

Set-TextValue "B50" "Monero"
Set-TextValue "C50" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D50" "143.90"
Set-TextValue "E50" "  +0.07%  "

Set-TextValue "D51" "2.56"
Set-TextValue "E51" "  +0.96%  "
